# This script applies the weekly update described in the commit:
# Two new weekly price rows for "Brócoli" at "Terminal Hortofrutícola Agro Chillán"
# are inserted at the top of the data table (row 518, pushing all the existing
# rows 518:626 down by two to 520:628), and the two newly inserted rows
# (518 and 519) are populated with this week's data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new blank rows before the current row 518; this pushes the
# existing rows 518-626 down to 520-628 (dimension becomes A1:R628).
$ws.Rows("518:519").Insert()

# --- New row 518: "Primera" quality ---
$ws.Range("A518").Value = 7
$ws.Range("B518").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C518").Value = "Ñuble"
$ws.Range("D518").Value = 45211
$ws.Range("E518").Value = 16
$ws.Range("F518").Value = 100112023
$ws.Range("G518").Value = "Brócoli"
$ws.Range("H518").Value = "Sin especificar"
$ws.Range("I518").Value = "Primera"
$ws.Range("J518").Value = 500
$ws.Range("K518").Value = 1200
$ws.Range("L518").Value = 1200
$ws.Range("M518").Value = 1200
$ws.Range("N518").Value = "$/unidad"
$ws.Range("O518").Value = "Región del Maule"
$ws.Range("P518").Value = 1200
$ws.Range("Q518").Value = 1
$ws.Range("R518").Value = "Hortaliza"

# --- New row 519: "Segunda" quality ---
$ws.Range("A519").Value = 7
$ws.Range("B519").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C519").Value = "Ñuble"
$ws.Range("D519").Value = 45211
$ws.Range("E519").Value = 16
$ws.Range("F519").Value = 100112023
$ws.Range("G519").Value = "Brócoli"
$ws.Range("H519").Value = "Sin especificar"
$ws.Range("I519").Value = "Segunda"
$ws.Range("J519").Value = 400
$ws.Range("K519").Value = 1000
$ws.Range("L519").Value = 1000
$ws.Range("M519").Value = 1000
$ws.Range("N519").Value = "$/unidad"
$ws.Range("O519").Value = "Región del Maule"
$ws.Range("P519").Value = 1000
$ws.Range("Q519").Value = 1
$ws.Range("R519").Value = "Hortaliza"
